$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.836.20"
$ws.Range("E2").Value = "  -0.61%  "

$ws.Range("D3").Value = "1.817.41"
$ws.Range("E3").Value = "  -0.73%  "

$ws.Range("D4").Value = "0.9929"
$ws.Range("E4").Value = "  -0.66%  "

$ws.Range("D5").Value = "242.11"
$ws.Range("E5").Value = "  +0.24%  "

$ws.Range("D6").Value = "0.6269"
$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("D7").Value = "0.9937"
$ws.Range("E7").Value = "  -0.70%  "

$ws.Range("D8").Value = "0.07440"
$ws.Range("E8").Value = "  -1.79%  "

$ws.Range("D9").Value = "0.2921"
$ws.Range("E9").Value = "  +0.48%  "

$ws.Range("D10").Value = "22.96"
$ws.Range("E10").Value = "  +1.22%  "

$ws.Range("D11").Value = "0.07656"
$ws.Range("E11").Value = "  -1.20%  "

$ws.Range("D12").Value = "1.821.19"
$ws.Range("E12").Value = "  -0.43%  "

$ws.Range("D13").Value = "4.974"
$ws.Range("E13").Value = "  +0.49%  "

$ws.Range("D14").Value = "0.6633"
$ws.Range("E14").Value = "  +0.32%  "

$ws.Range("D15").Value = "82.56"
$ws.Range("E15").Value = "  +0.24%  "

$ws.Range("D16").Value = "0.000009624"
$ws.Range("E16").Value = "  +1.65%  "

$ws.Range("D17").Value = "5.994"
$ws.Range("E17").Value = "  +0.75%  "

$ws.Range("D18").Value = "28.873.76"
$ws.Range("E18").Value = "  -0.46%  "

$ws.Range("D19").Value = "12.50"
$ws.Range("E19").Value = "  +1.60%  "

$ws.Range("D20").Value = "223.16"
$ws.Range("E20").Value = "  -0.27%  "

$ws.Range("D21").Value = "0.9934"
$ws.Range("E21").Value = "  -0.69%  "

$ws.Range("D22").Value = "7.085"
$ws.Range("E22").Value = "  -1.45%  "

$ws.Range("D23").Value = "0.9944"
$ws.Range("E23").Value = "  -0.66%  "

$ws.Range("D24").Value = "158.66"
$ws.Range("E24").Value = "  -0.62%  "

$ws.Range("D25").Value = "0.1409"
$ws.Range("E25").Value = "  +3.59%  "

$ws.Range("D26").Value = "8.436"
$ws.Range("E26").Value = "  +0.32%  "

$ws.Range("D27").Value = "17.80"
$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("D28").Value = "1.492"
$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("D29").Value = "4.098"
$ws.Range("E29").Value = "  +1.10%  "

$ws.Range("D30").Value = "4.031"
$ws.Range("E30").Value = "  +0.32%  "

$ws.Range("D31").Value = "0.05434"
$ws.Range("E31").Value = "  +4.88%  "

$ws.Range("D32").Value = "1.192"
$ws.Range("E32").Value = "  -0.47%  "

$ws.Range("D33").Value = "1.845"
$ws.Range("E33").Value = "  -0.06%  "

$ws.Range("D34").Value = "0.7382"
$ws.Range("E34").Value = "  +0.32%  "

$ws.Range("D35").Value = "1.128"
$ws.Range("E35").Value = "  -1.48%  "

$ws.Range("D36").Value = "2.601"
$ws.Range("E36").Value = "  -3.63%  "

$ws.Range("D37").Value = "1.231.79"
$ws.Range("E37").Value = "  -2.20%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "2.734"
$ws.Range("E38").Value = "  -0.85%  "

$ws.Range("D39").Value = "0.01775"
$ws.Range("E39").Value = "  -0.50%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "6.672"
$ws.Range("E40").Value = "  +6.92%  "

$ws.Range("D41").Value = "0.8937"
$ws.Range("E41").Value = "  +0.50%  "

$ws.Range("D42").Value = "0.9945"
$ws.Range("E42").Value = "  -0.62%  "

$ws.Range("D43").Value = "100.84"
$ws.Range("E43").Value = "  -0.62%  "

$ws.Range("D44").Value = "1.973.84"
$ws.Range("E44").Value = "  -0.16%  "

$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.00000000123"
$ws.Range("E45").Value = "  +2.50%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "64.94"
$ws.Range("E46").Value = "  +0.95%  "

$ws.Range("D47").Value = "0.5049"
$ws.Range("E47").Value = "  -1.21%  "

$ws.Range("D48").Value = "0.4021"
$ws.Range("E48").Value = "  +1.16%  "

$ws.Range("D49").Value = "8.934"
$ws.Range("E49").Value = "  +0.90%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "1.663"
$ws.Range("E50").Value = "  +2.19%  "

$ws.Range("B51").Value = "XinFinNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D51").Value = "0.07177"
$ws.Range("E51").Value = "  +1.65%  "
